$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from just before "Modify the scoring
#    algorithm..." to the very start of the document (before the
#    first run of the title paragraph "Powerup Coding Challenges").
# -----------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Adding a bookmark directly at document position 0 with a collapsed
# range causes it to balloon to cover the whole first paragraph, so
# instead: temporarily insert an extra empty placeholder paragraph at
# the very start, put a placeholder character in it, anchor the
# (still-collapsed) bookmark right after that placeholder paragraph,
# then delete the whole placeholder paragraph again. Bookmark gravity
# keeps the bookmark pinned to position 0 once the placeholder is gone.
$startRange = $d.Range(0, 0)
$startRange.InsertParagraphBefore()

$placeholderCharRange = $d.Range(0, 0)
$placeholderCharRange.InsertBefore("X")

$afterPlaceholder = $d.Paragraphs(2).Range.Start
$bmRange = $d.Range($afterPlaceholder, $afterPlaceholder)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholderRange = $d.Paragraphs(1).Range
$placeholderRange.Delete()

# -----------------------------------------------------------------
# 2) Merge the two runs "New Programmer" + "s" into a single run
#    "New Programmers" inside the table cell.
# -----------------------------------------------------------------
$table = $d.Tables(1)
$cellRange = $table.Cell(2, 1).Range
$cellRange.Find.Execute("New Programmers", $true, $true, $false, $false, $false, `
    $true, 1, $false, "New Programmers", 2)

# -----------------------------------------------------------------
# 3) Update the footer DATE field's cached result text.
# -----------------------------------------------------------------
$footer = $d.Sections(1).Footers(1)
$footerRange = $footer.Range
$footerRange.Find.Execute("6/12/2018 5:59 AM", $false, $false, $false, $false, $false, `
    $true, 1, $false, "6/16/2018 2:09 PM", 2)
